# Translate test workbook contents from Norwegian to English, and refresh
# the cell/row formatting that Google Sheets applied on round-trip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3) ---
$ws.Range("A3").Value = "Debit account"
$ws.Range("B3").Value = "Credit account"
$ws.Range("C3").Value = "Amount"
$ws.Range("D3").Value = "Date"
$ws.Range("E3").Value = "Text"

# --- Merged-cell labels ---
$ws.Range("A4").Value = "Short merged"
$ws.Range("A5").Value = "Long merged cell spanning two columns"

# --- Description column (E) ---
$ws.Range("E4").Value = "I am a description"
$ws.Range("E5").Value = "Had to move some money"
$ws.Range("E6").Value = "Had to move some money"
$ws.Range("E7").Value = "Had to move some money"
$ws.Range("E8").Value = "Hello"
$ws.Range("E9").Value = "Had to move some money"
$ws.Range("E10").Value = "Had to move some money"
$ws.Range("E11").Value = "Had to move some money"
$ws.Range("E12").Value = "Had to move some money"
$ws.Range("E13").Value = "Had to move some money"
$ws.Range("E14").Value = "Had to move some money"
$ws.Range("E15").Value = "Had to move some money"
$ws.Range("E16").Value = "Had to move some money"
